# Added effort estimates to historical work:
#  - renamed "Sheet1" -> "Estimates"
#  - added a new "Effort" sheet after it with nvessels/sets/obs data
#  - bolded the Effort header row (a fresh bold+automatic-color font/style)

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Estimates"

# Insert the new sheet right after "Estimates" so tab order is Estimates, Effort.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Effort"

# Header row
$ws2.Range("A1").Value = "year"
$ws2.Range("B1").Value = "nvessels"
$ws2.Range("C1").Value = "sets_tot_est"
$ws2.Range("D1").Value = "sets_obs"
$ws2.Range("E1").Value = "obs_perc"
$ws2.Range("F1").Value = "notes"
[void]($ws2.Range("A1:F1").Font.Bold = $true)

# Data row
$ws2.Range("A2").Value = 2007
$ws2.Range("B2").Value = 58
$ws2.Range("C2").Value = 1387
$ws2.Range("D2").Value = 248
$ws2.Range("E2").Value = 17.8
$ws2.Range("F2").Value = "effort based on 2006 logbook data"

# Size columns to fit the new content, like the author did in Excel.
$ws2.Range("A1:F2").EntireColumn.AutoFit()

# Leave the selection where the author's last edit left it.
[void]$ws2.Range("F3").Select()
